# Auto-generated Excel COM-interop script applying Chocobo_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 52
$ws.Range("I11").Value = 52
$ws.Range("K11").Value = 52
$ws.Range("M11").Value = 88
$ws.Range("H41").Value = 187.27272
$ws.Range("J41").Value = 181.8
$ws.Range("L41").Value = 181.8
$ws.Range("N41").Value = -1061.8
$ws.Range("H53").Value = 583.35
$ws.Range("I53").Value = 236.9
$ws.Range("J53").Value = 929.8
$ws.Range("K53").Value = 236.9
$ws.Range("L53").Value = 929.8
$ws.Range("M53").Value = 400.1
$ws.Range("N53").Value = -2203.8
$ws.Range("H93").Value = 24489
$ws.Range("J93").Value = 24489
$ws.Range("L93").Value = 24489
$ws.Range("N93").Value = -29481
$ws.Range("H103").Value = 5048.9
$ws.Range("I103").Value = 796.7692
$ws.Range("J103").Value = 12945.714
$ws.Range("K103").Value = 2390.3076
$ws.Range("L103").Value = 38837.142
$ws.Range("M103").Value = -1804.3076
$ws.Range("N103").Value = -40009.142
$ws.Range("H112").Value = 1310.7384
$ws.Range("J112").Value = 1310.7384
$ws.Range("L112").Value = 3932.2152
$ws.Range("N112").Value = -6148.2152
$ws.Range("H129").Value = 1766.6123
$ws.Range("J129").Value = 1965.279
$ws.Range("L129").Value = 5895.837
$ws.Range("N129").Value = -15895.837
$ws.Range("H131").Value = 5856.0835
$ws.Range("I131").Value = 3028.3
$ws.Range("K131").Value = 9084.900000000001
$ws.Range("M131").Value = -4044.900000000001
$ws.Range("H137").Value = 783594.5
$ws.Range("I137").Value = 1644954.4
$ws.Range("J137").Value = 2987.1875
$ws.Range("K137").Value = 4934863.199999999
$ws.Range("L137").Value = 8961.5625
$ws.Range("M137").Value = -4932313.199999999
$ws.Range("N137").Value = -14061.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4403.443
$ws.Range("I32").Value = 4864.396
$ws.Range("J32").Value = 3397.7273
$ws.Range("K32").Value = 4864.396
$ws.Range("L32").Value = 3397.7273
$ws.Range("M32").Value = -4577.396
$ws.Range("N32").Value = -3971.7273
$ws.Range("H94").Value = 50000
$ws.Range("J94").Value = 50000
$ws.Range("L94").Value = 50000
$ws.Range("N94").Value = -51802

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 833.2
$ws.Range("I94").Value = 783.5
$ws.Range("J94").Value = 866.3333
$ws.Range("K94").Value = 783.5
$ws.Range("L94").Value = 866.3333
$ws.Range("M94").Value = -332.5
$ws.Range("N94").Value = -1768.3333
$ws.Range("H122").Value = 42996
$ws.Range("J122").Value = 42996
$ws.Range("L122").Value = 42996
$ws.Range("N122").Value = -52796

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 179696.88
$ws.Range("I31").Value = 466767.06
$ws.Range("K31").Value = 466767.06
$ws.Range("M31").Value = -466472.06
$ws.Range("H34").Value = 179696.88
$ws.Range("I34").Value = 466767.06
$ws.Range("K34").Value = 466767.06
$ws.Range("M34").Value = -466565.06
$ws.Range("H105").Value = 1755.9231
$ws.Range("I105").Value = 1536.3334
$ws.Range("J105").Value = 2250
$ws.Range("K105").Value = 1536.3334
$ws.Range("L105").Value = 2250
$ws.Range("M105").Value = 210.6666
$ws.Range("N105").Value = -5744
$ws.Range("H132").Value = 3457
$ws.Range("I132").Value = 1640
$ws.Range("J132").Value = 7999.5
$ws.Range("K132").Value = 4920
$ws.Range("L132").Value = 23998.5
$ws.Range("M132").Value = -2390
$ws.Range("N132").Value = -29058.5
$ws.Range("H134").Value = 1521.697
$ws.Range("I134").Value = 1014.9474
$ws.Range("J134").Value = 2209.4285
$ws.Range("K134").Value = 3044.8422
$ws.Range("L134").Value = 6628.2855
$ws.Range("M134").Value = -509.8422
$ws.Range("N134").Value = -11698.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 79.73333
$ws.Range("J12").Value = 127.666664
$ws.Range("L12").Value = 382.999992
$ws.Range("N12").Value = -728.999992
$ws.Range("H113").Value = 445.86792
$ws.Range("I113").Value = 463.53845
$ws.Range("J113").Value = 428.85184
$ws.Range("K113").Value = 1390.61535
$ws.Range("L113").Value = 1286.55552
$ws.Range("M113").Value = 779.38465
$ws.Range("N113").Value = -5626.55552
$ws.Range("H129").Value = 1460.88
$ws.Range("J129").Value = 1898.4546
$ws.Range("L129").Value = 5695.3638
$ws.Range("N129").Value = -15695.3638
$ws.Range("H131").Value = 5000885
$ws.Range("I131").Value = 62500590
$ws.Range("J131").Value = 910.25
$ws.Range("K131").Value = 187501770
$ws.Range("L131").Value = 2730.75
$ws.Range("M131").Value = -187496730
$ws.Range("N131").Value = -12810.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 8009.091
$ws.Range("I122").Value = 3300
$ws.Range("J122").Value = 11933.333
$ws.Range("K122").Value = 9900
$ws.Range("L122").Value = 35799.999
$ws.Range("M122").Value = -7450
$ws.Range("N122").Value = -40699.999
$ws.Range("H132").Value = 3844.5789
$ws.Range("I132").Value = 2668.6086
$ws.Range("J132").Value = 5647.7334
$ws.Range("K132").Value = 8005.825800000001
$ws.Range("L132").Value = 16943.2002
$ws.Range("M132").Value = -5475.825800000001
$ws.Range("N132").Value = -22003.2002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3634.3333
$ws.Range("I7").Value = 1984
$ws.Range("J7").Value = 6935
$ws.Range("K7").Value = 1984
$ws.Range("L7").Value = 6935
$ws.Range("M7").Value = -1872
$ws.Range("N7").Value = -7159
$ws.Range("H22").Value = 2348.8125
$ws.Range("I22").Value = 1575.8
$ws.Range("J22").Value = 2700.182
$ws.Range("K22").Value = 1575.8
$ws.Range("L22").Value = 2700.182
$ws.Range("M22").Value = -1280.8
$ws.Range("N22").Value = -3290.182
$ws.Range("H27").Value = 2348.8125
$ws.Range("I27").Value = 1575.8
$ws.Range("J27").Value = 2700.182
$ws.Range("K27").Value = 1575.8
$ws.Range("L27").Value = 2700.182
$ws.Range("M27").Value = -1468.8
$ws.Range("N27").Value = -2914.182
$ws.Range("H40").Value = 5052
$ws.Range("I40").Value = 4620.7334
$ws.Range("K40").Value = 4620.7334
$ws.Range("M40").Value = -4484.7334
$ws.Range("H55").Value = 344
$ws.Range("I55").Value = 292.16666
$ws.Range("K55").Value = 292.16666
$ws.Range("M55").Value = -119.16666
$ws.Range("H68").Value = 992.81177
$ws.Range("I68").Value = 912.481
$ws.Range("J68").Value = 2050.5
$ws.Range("K68").Value = 912.481
$ws.Range("L68").Value = 2050.5
$ws.Range("M68").Value = -163.481
$ws.Range("N68").Value = -3548.5
$ws.Range("H71").Value = 992.81177
$ws.Range("I71").Value = 912.481
$ws.Range("J71").Value = 2050.5
$ws.Range("K71").Value = 4562.405
$ws.Range("L71").Value = 10252.5
$ws.Range("M71").Value = -818.4049999999997
$ws.Range("N71").Value = -17740.5
$ws.Range("H122").Value = 6065.8887
$ws.Range("I122").Value = 3898.6
$ws.Range("J122").Value = 8775
$ws.Range("K122").Value = 11695.8
$ws.Range("L122").Value = 26325
$ws.Range("M122").Value = -9245.799999999999
$ws.Range("N122").Value = -31225
$ws.Range("H126").Value = 3634.3333
$ws.Range("I126").Value = 1984
$ws.Range("J126").Value = 6935
$ws.Range("K126").Value = 5952
$ws.Range("L126").Value = 20805
$ws.Range("M126").Value = -3482
$ws.Range("N126").Value = -25745
$ws.Range("H132").Value = 3842.423
$ws.Range("I132").Value = 2806.7058
$ws.Range("J132").Value = 5798.778
$ws.Range("K132").Value = 8420.117400000001
$ws.Range("L132").Value = 17396.334
$ws.Range("M132").Value = -5890.117400000001
$ws.Range("N132").Value = -22456.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6745.769
$ws.Range("I122").Value = 3100
$ws.Range("J122").Value = 7408.636
$ws.Range("K122").Value = 9300
$ws.Range("L122").Value = 22225.908
$ws.Range("M122").Value = -6850
$ws.Range("N122").Value = -27125.908
$ws.Range("H126").Value = 820938.75
$ws.Range("I126").Value = 1956
$ws.Range("K126").Value = 5868
$ws.Range("M126").Value = -3398
$ws.Range("H132").Value = 9014423
$ws.Range("I132").Value = 8684.929
$ws.Range("J132").Value = 14496176
$ws.Range("K132").Value = 26054.787
$ws.Range("L132").Value = 43488528
$ws.Range("M132").Value = -23524.787
$ws.Range("N132").Value = -43493588
